# Update "想去人数" (F column) counts on both the "展览" sheet and the
# aggregated "全部类型" sheet, per the generated-data refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows 2-9, column F)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 699
$ws1.Range("F3").Value = 33
$ws1.Range("F4").Value = 236
$ws1.Range("F5").Value = 2264
$ws1.Range("F6").Value = 48
$ws1.Range("F7").Value = 3480
$ws1.Range("F8").Value = 463
$ws1.Range("F9").Value = 872

# Sheet "全部类型" (rows 2-10, column F; row 4 unchanged)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 699
$ws4.Range("F3").Value = 33
$ws4.Range("F5").Value = 236
$ws4.Range("F6").Value = 2264
$ws4.Range("F7").Value = 48
$ws4.Range("F8").Value = 3480
$ws4.Range("F9").Value = 463
$ws4.Range("F10").Value = 872

$wb.Save()
